# Auto-committed on 2022/05/10 週二
#
# DBS sheet gains a new lookup-definition row ("maxLogNoFirst" / "CustNo >" /
# "LogNo DESC", reusing the existing "LogNo DESC" entry from row 6), and the
# two sheets' saved cursor/selection positions move to where the user last
# clicked (DBD -> D20, DBS -> B7).

$wb    = $excel.ActiveWorkbook
$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# Append row 7 on DBS, right under the existing findCustNoFirst row (row 6).
# Leaving style alone lets the new cells inherit the sheet's column style
# (same s="1" as every other data row) instead of minting a new style.
$wsDBS.Range("A7").Value = "maxLogNoFirst"
$wsDBS.Range("B7").Value = "CustNo >"
$wsDBS.Range("C7").Value = $wsDBS.Range("C6").Value()

# Restore the cell selections recorded in the saved workbook.
[void]$wsDBD.Range("D20").Select()
[void]$wsDBS.Range("B7").Select()
